$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$m.Shapes.Item(1).TextFrame.TextRange.Text = $m.Shapes.Item(1).TextFrame.TextRange.Text
